# Apply scheduled-runner price/profit updates to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 649.8378
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 646.97144
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 1940.91432
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -2276.91432

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 5229.4473
$ws.Range("I113").Value = 7551.706
$ws.Range("J113").Value = 3349.524
$ws.Range("K113").Value = 7551.706
$ws.Range("L113").Value = 3349.524
$ws.Range("M113").Value = -4297.706
$ws.Range("N113").Value = -9857.523999999999

$ws.Range("H127").Value = 2541.8235
$ws.Range("I127").Value = 2166.6667
$ws.Range("J127").Value = 2622.2144
$ws.Range("K127").Value = 6500.000100000001
$ws.Range("L127").Value = 7866.6432
$ws.Range("M127").Value = -1540.000100000001
$ws.Range("N127").Value = -17786.6432

$ws.Range("H129").Value = 875.2353000000001
$ws.Range("I129").Value = 611.4286
$ws.Range("J129").Value = 1059.9
$ws.Range("K129").Value = 1834.2858
$ws.Range("L129").Value = 3179.7
$ws.Range("M129").Value = 3165.7142
$ws.Range("N129").Value = -13179.7

$ws.Range("H137").Value = 271981.72
$ws.Range("I137").Value = 701274.6
$ws.Range("J137").Value = 1686.1482
$ws.Range("K137").Value = 2103823.8
$ws.Range("L137").Value = 5058.444600000001
$ws.Range("M137").Value = -2101273.8
$ws.Range("N137").Value = -10158.4446

$ws.Range("H138").Value = 1788.98
$ws.Range("I138").Value = 795.0278
$ws.Range("J138").Value = 2348.0781
$ws.Range("K138").Value = 2385.0834
$ws.Range("L138").Value = 7044.2343
$ws.Range("M138").Value = 2754.9166
$ws.Range("N138").Value = -17324.2343

$ws.Range("H141").Value = 7759.3335
$ws.Range("I141").Value = 8491.538
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 25474.614
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -20294.614
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5698.096
$ws.Range("I32").Value = 5145.3696
$ws.Range("J32").Value = 9935.666999999999
$ws.Range("K32").Value = 5145.3696
$ws.Range("L32").Value = 9935.666999999999
$ws.Range("M32").Value = -4858.3696
$ws.Range("N32").Value = -10509.667

$ws.Range("H74").Value = 3305.8975
$ws.Range("I74").Value = 366.84616
$ws.Range("J74").Value = 9184
$ws.Range("K74").Value = 366.84616
$ws.Range("L74").Value = 9184
$ws.Range("M74").Value = 507.15384
$ws.Range("N74").Value = -10932

$ws.Range("H77").Value = 3305.8975
$ws.Range("I77").Value = 366.84616
$ws.Range("J77").Value = 9184
$ws.Range("K77").Value = 1834.2308
$ws.Range("L77").Value = 45920
$ws.Range("M77").Value = 2533.7692
$ws.Range("N77").Value = -54656

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 43480396
$ws.Range("I107").Value = 100002104
$ws.Range("J107").Value = 2159.4614
$ws.Range("K107").Value = 100002104
$ws.Range("L107").Value = 2159.4614
$ws.Range("M107").Value = -100000184
$ws.Range("N107").Value = -5999.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9337.529
$ws.Range("I31").Value = 1150.6957
$ws.Range("J31").Value = 16062.429
$ws.Range("K31").Value = 1150.6957
$ws.Range("L31").Value = 16062.429
$ws.Range("M31").Value = -855.6957
$ws.Range("N31").Value = -16652.429

$ws.Range("H34").Value = 9337.529
$ws.Range("I34").Value = 1150.6957
$ws.Range("J34").Value = 16062.429
$ws.Range("K34").Value = 1150.6957
$ws.Range("L34").Value = 16062.429
$ws.Range("M34").Value = -948.6957
$ws.Range("N34").Value = -16466.429

$ws.Range("H58").Value = 3348489.5
$ws.Range("I58").Value = 3888606.8
$ws.Range("J58").Value = 17766.334
$ws.Range("K58").Value = 3888606.8
$ws.Range("L58").Value = 17766.334
$ws.Range("M58").Value = -3888403.8
$ws.Range("N58").Value = -18172.334

$ws.Range("H136").Value = 3348489.5
$ws.Range("I136").Value = 3888606.8
$ws.Range("J136").Value = 17766.334
$ws.Range("K136").Value = 11665820.4
$ws.Range("L136").Value = 53299.00199999999
$ws.Range("M136").Value = -11663270.4
$ws.Range("N136").Value = -58399.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15288539
$ws.Range("I131").Value = 90909496
$ws.Range("J131").Value = 1424697.8
$ws.Range("K131").Value = 272728488
$ws.Range("L131").Value = 4274093.4
$ws.Range("M131").Value = -272723448
$ws.Range("N131").Value = -4284173.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 35715108
$ws.Range("I122").Value = 52632292
$ws.Range("J122").Value = 1054.1111
$ws.Range("K122").Value = 157896876
$ws.Range("L122").Value = 3162.3333
$ws.Range("M122").Value = -157894426
$ws.Range("N122").Value = -8062.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1823.2084
$ws.Range("I7").Value = 1121.1
$ws.Range("J7").Value = 2324.7144
$ws.Range("K7").Value = 1121.1
$ws.Range("L7").Value = 2324.7144
$ws.Range("M7").Value = -1009.1
$ws.Range("N7").Value = -2548.7144

$ws.Range("H61").Value = 3314.1333
$ws.Range("I61").Value = 2500.8
$ws.Range("J61").Value = 3720.8
$ws.Range("K61").Value = 2500.8
$ws.Range("L61").Value = 3720.8
$ws.Range("M61").Value = -2298.8
$ws.Range("N61").Value = -4124.8

$ws.Range("H113").Value = 3314.1333
$ws.Range("I113").Value = 2500.8
$ws.Range("J113").Value = 3720.8
$ws.Range("K113").Value = 2500.8
$ws.Range("L113").Value = 3720.8
$ws.Range("M113").Value = -330.8000000000002
$ws.Range("N113").Value = -8060.8

$ws.Range("H122").Value = 4042
$ws.Range("I122").Value = 1504
$ws.Range("J122").Value = 4465
$ws.Range("K122").Value = 4512
$ws.Range("L122").Value = 13395
$ws.Range("M122").Value = -2062
$ws.Range("N122").Value = -18295

$ws.Range("H126").Value = 1823.2084
$ws.Range("I126").Value = 1121.1
$ws.Range("J126").Value = 2324.7144
$ws.Range("K126").Value = 3363.3
$ws.Range("L126").Value = 6974.1432
$ws.Range("M126").Value = -893.2999999999997
$ws.Range("N126").Value = -11914.1432

$ws.Range("H132").Value = 4547012.5
$ws.Range("I132").Value = 6061532
$ws.Range("J132").Value = 3453.6365
$ws.Range("K132").Value = 18184596
$ws.Range("L132").Value = 10360.9095
$ws.Range("M132").Value = -18182066
$ws.Range("N132").Value = -15420.9095

$ws.Range("H136").Value = 4000.2917
$ws.Range("I136").Value = 4962.394
$ws.Range("J136").Value = 1883.6666
$ws.Range("K136").Value = 14887.182
$ws.Range("L136").Value = 5650.9998
$ws.Range("M136").Value = -12337.182
$ws.Range("N136").Value = -10750.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 89286616
$ws.Range("I126").Value = 113637220
$ws.Range("J126").Value = 1066.6666
$ws.Range("K126").Value = 340911660
$ws.Range("L126").Value = 3199.9998
$ws.Range("M126").Value = -340909190
$ws.Range("N126").Value = -8139.9998

$ws.Range("H136").Value = 11904786
$ws.Range("I136").Value = 6052773
$ws.Range("J136").Value = 38463924
$ws.Range("K136").Value = 18158319
$ws.Range("L136").Value = 115391772
$ws.Range("M136").Value = -18155769
$ws.Range("N136").Value = -115396872
